$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($refAddr, $val) {
    $c = $ws.Range($refAddr)
    $c.Formula = '="' + $val + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" "30.551.01"
Set-TextValue "E2" "  -0.13%  "
Set-TextValue "D3" "1.873.67"
Set-TextValue "E3" "  -0.94%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "247.75"
Set-TextValue "E5" "  +1.15%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  +0.06%  "
Set-TextValue "D7" "0.4730"
Set-TextValue "E7" "  -1.07%  "
Set-TextValue "D8" "0.2897"
Set-TextValue "E8" "  -0.17%  "
Set-TextValue "D9" "0.06474"
Set-TextValue "E9" "  -1.27%  "
Set-TextValue "E10" "  +2.80%  "
Set-TextValue "D11" "0.07725"
Set-TextValue "E11" "  -0.75%  "
Set-TextValue "D12" "0.7412"
Set-TextValue "E12" "  +0.10%  "
Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.873.12"
Set-TextValue "E13" "  -0.93%  "
Set-TextValue "B14" "Litecoin"
Set-TextValue "C14" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D14" "96.00"
Set-TextValue "E14" "  -0.33%  "
Set-TextValue "D15" "5.172"
Set-TextValue "D16" "274.18"
Set-TextValue "E16" "  -1.30%  "
Set-TextValue "D17" "30.612.35"
Set-TextValue "E17" "  +0.12%  "
Set-TextValue "E18" "  -3.07%  "
Set-TextValue "E19" "  +0.08%  "
Set-TextValue "D20" "0.000007470"
Set-TextValue "E20" "  -2.08%  "
Set-TextValue "D21" "2.120.28"
Set-TextValue "E21" "  -0.97%  "
Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  +0.04%  "
Set-TextValue "D23" "5.206"
Set-TextValue "E23" "  -2.06%  "
Set-TextValue "D24" "6.168"
Set-TextValue "E24" "  -0.96%  "
Set-TextValue "B25" "Cosmos"
Set-TextValue "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D25" "9.186"
Set-TextValue "E25" "  -1.51%  "
Set-TextValue "B26" "Monero"
Set-TextValue "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D26" "165.14"
Set-TextValue "E26" "  -0.11%  "
Set-TextValue "D27" "18.65"
Set-TextValue "E27" "  -2.31%  "
Set-TextValue "D28" "1.902"
Set-TextValue "E28" "  -4.77%  "
Set-TextValue "D29" "0.09956"
Set-TextValue "E29" "  -0.32%  "
Set-TextValue "E30" "  -2.51%  "
Set-TextValue "D31" "1.508"
Set-TextValue "E31" "  -0.50%  "
Set-TextValue "D32" "4.233"
Set-TextValue "E32" "  -2.75%  "
Set-TextValue "D33" "4.085"
Set-TextValue "E33" "  -0.87%  "
Set-TextValue "E34" "  -0.30%  "
Set-TextValue "D35" "1.118"
Set-TextValue "E35" "  -1.39%  "
Set-TextValue "E36" "  -1.94%  "
Set-TextValue "E37" "  -0.04%  "
Set-TextValue "D38" "0.01846"
Set-TextValue "E38" "  -1.09%  "
Set-TextValue "D39" "2.754"
Set-TextValue "E39" "  -0.25%  "
Set-TextValue "D40" "6.262"
Set-TextValue "E40" "  -4.15%  "
Set-TextValue "D41" "73.24"
Set-TextValue "E41" "  +3.31%  "
Set-TextValue "D42" "1.969"
Set-TextValue "E42" "  +2.02%  "
Set-TextValue "E43" "  +0.06%  "
Set-TextValue "D44" "0.4157"
Set-TextValue "E44" "  -1.13%  "
Set-TextValue "D45" "0.8334"
Set-TextValue "E45" "  -1.97%  "
Set-TextValue "D46" "101.14"
Set-TextValue "E46" "  -1.70%  "
Set-TextValue "D47" "9.330"
Set-TextValue "B48" "Aptos"
Set-TextValue "C48" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "6.975"
Set-TextValue "E48" "  -2.57%  "
Set-TextValue "B49" "Elrond"
Set-TextValue "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D49" "35.29"
Set-TextValue "E49" "  -0.52%  "
Set-TextValue "D50" "912.31"
Set-TextValue "E50" "  -1.89%  "
Set-TextValue "D51" "0.05663"
Set-TextValue "E51" "  +1.04%  "

$excel.CutCopyMode = 0
